$d = $word.ActiveDocument

# Locate the empty paragraph that sits right after the (last) table:
#   <w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/>
#        <w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p>
# It has no run content (just the paragraph mark), 13.8pt auto line spacing
# (w:line="276" w:lineRule="auto" -> 276/20 = 13.8), and both-justified alignment.
# NOTE: deliberately avoid touching $d.Tables here - doing so before walking
# $d.Paragraphs throws off Range.Start/End bookkeeping for later paragraphs.

$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt.Length -le 2) {
        if ($p.Format.LineSpacing -eq 13.8 -and $p.Format.Alignment -eq 3) {
            $target = $p
        }
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
